$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.356678009033203
$ws.Range("B1").Value = 2.639779806137085
$ws.Range("C1").Value = 2.681990146636963
$ws.Range("D1").Value = 3.404010534286499
$ws.Range("E1").Value = 1.93707799911499
